# Changes in Work Order/System Setup
# The "Create Disassembly WO" sheet had its second data row (row 2, the
# "MS-Disassembly (NO Track)" component) removed. Deleting the entire row
# shifts the remaining rows up, the dimension shrinks by one row, and the
# active selection moves to the new row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create Disassembly WO")

# Remove the obsolete "MS-Disassembly (NO Track)" row - this shifts
# row 3 -> row 2 and row 4 -> row 3 automatically.
$ws.Rows.Item(2).Delete()

# Reflect the new selection state left behind after the row delete.
$ws.Range("A2:XFD2").Select()
